$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 747
$ws.Range("F3").Value = 14120
$ws.Range("F4").Value = 14120
$ws.Range("F5").Value = 14180
$ws.Range("F7").Value = 1392
$ws.Range("F8").Value = 5854
$ws.Range("F9").Value = 980
$ws.Range("F14").Value = 1529
$ws.Range("F15").Value = 431
$ws.Range("F17").Value = 1186
$ws.Range("F18").Value = 1806
$ws.Range("G18").Value = 100
$ws.Range("F19").Value = 912
$ws.Range("F20").Value = 34
$ws.Range("F21").Value = 2265
$ws.Range("F23").Value = 803
$ws.Range("F24").Value = 3300
$ws.Range("F26").Value = 307
$ws.Range("F27").Value = 2363
$ws.Range("F28").Value = 583
$ws.Range("F29").Value = 117
$ws.Range("F30").Value = 1341
$ws.Range("F31").Value = 1776
$ws.Range("F32").Value = 1068
$ws.Range("F33").Value = 1366
$ws.Range("F35").Value = 140
$ws.Range("F36").Value = 4739
$ws.Range("F37").Value = 4796
$ws.Range("F38").Value = 298
$ws.Range("F39").Value = 156
$ws.Range("F41").Value = 682
$ws.Range("F42").Value = 3279
$ws.Range("F46").Value = 99
$ws.Range("F47").Value = 68
$ws.Range("F48").Value = 4415
$ws.Range("F49").Value = 567
$ws.Range("F50").Value = 281

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 46
$ws.Range("F4").Value = 117
$ws.Range("F14").Value = 13

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 7500
$ws.Range("F3").Value = 230
$ws.Range("F4").Value = 736

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 7500
$ws.Range("F3").Value = 747
$ws.Range("F4").Value = 230
$ws.Range("F5").Value = 736
$ws.Range("F6").Value = 14120
$ws.Range("F7").Value = 14180
$ws.Range("F9").Value = 1392
$ws.Range("F10").Value = 5854
$ws.Range("F11").Value = 980
$ws.Range("F12").Value = 117
$ws.Range("F15").Value = 1529
$ws.Range("F16").Value = 432
$ws.Range("F17").Value = 1186
$ws.Range("F18").Value = 1806
$ws.Range("G18").Value = 100
$ws.Range("F19").Value = 912
$ws.Range("F20").Value = 34
$ws.Range("F21").Value = 3300
$ws.Range("F22").Value = 307
$ws.Range("F23").Value = 2363
$ws.Range("F24").Value = 583
$ws.Range("F25").Value = 117
$ws.Range("F27").Value = 1776
$ws.Range("F31").Value = 1068
$ws.Range("F32").Value = 1367
$ws.Range("F34").Value = 4739
$ws.Range("F35").Value = 4796
$ws.Range("F36").Value = 298
$ws.Range("F37").Value = 156
$ws.Range("F39").Value = 682
$ws.Range("F40").Value = 3279
$ws.Range("F43").Value = 99
$ws.Range("F45").Value = 68
$ws.Range("F46").Value = 4415
$ws.Range("F47").Value = 567
$ws.Range("F48").Value = 281
